# Generate Report for Handoff
# Replaces the two previously-handed-back source files with a newly
# generated handoff report: new file names/UUIDs, new commit hash for the
# handoff packages, "Ready for handoff" status, fresh handoff timestamps,
# and handback info reset back to "not yet handed back".

$wb = $excel.ActiveWorkbook

$oldMd1 = "760568b4-a7e0-499e-a287-a54bd954753e.md"
$oldMd2 = "bccb55ee-d540-40a4-837e-fac60252d379.md"
$newMd1 = "e04a53de-7ccd-4d91-8db5-1ad4f9b4db22.md"
$newMd2 = "ffff345ff04a-2c00-4709-9758-4445ee9587fe.md"

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

$notHandedBack = "0001-01-01 00:00:00"

function Set-RowHyperlink($ws, $cellAddr, $newDisplay, $newUrl) {
    foreach ($hl in @($ws.Hyperlinks)) {
        if ($hl.Range.Address() -eq $cellAddr) {
            $hl.TextToDisplay = $newDisplay
            $hl.Address = $newUrl
        }
    }
}

function Remove-RowHyperlink($ws, $cellAddr) {
    foreach ($hl in @($ws.Hyperlinks)) {
        if ($hl.Range.Address() -eq $cellAddr) {
            $hl.Delete()
        }
    }
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
Set-RowHyperlink $wsOverview '$A$2' $newMd1 "https://github.com/OpenLocalizationTest/oltest/blob/7cae2979785f791e618474df3678176f73d707b9/e2e/$newMd1"

$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
Set-RowHyperlink $wsOverview '$A$3' $newMd2 "https://github.com/OpenLocalizationTest/oltest/blob/7cae2979785f791e618474df3678176f73d707b9/e2e/$newMd2"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$newXlf1Zh = "e04a53de-7ccd-4d91-8db5-1ad4f9b4db22.ce213699e41debf6b1b6a57a302349c00d6d9f03.zh-cn.xlf"
$handoffTimeZh = "2016-03-08 17:06:34"

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("C2").Value = $newXlf1Zh
$wsZh.Range("D2").Value = $handoffTimeZh
$wsZh.Range("G2").Value = $notHandedBack
Set-RowHyperlink $wsZh '$A$2' $newMd1 "https://github.com/OpenLocalizationTest/oltest/blob/7cae2979785f791e618474df3678176f73d707b9/e2e/$newMd1"
Set-RowHyperlink $wsZh '$C$2' $newXlf1Zh "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70a010d6535a8f2af430b4e65d2a91828e4ac020/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlf1Zh"
Remove-RowHyperlink $wsZh '$E$2'
Remove-RowHyperlink $wsZh '$F$2'
$wsZh.Range("E2").Clear()
$wsZh.Range("F2").Clear()

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("B3").Value = $newStatus
$wsZh.Range("C3").Value = $newXlf1Zh
$wsZh.Range("D3").Value = $handoffTimeZh
$wsZh.Range("G3").Value = $notHandedBack
Set-RowHyperlink $wsZh '$A$3' $newMd2 "https://github.com/OpenLocalizationTest/oltest/blob/7cae2979785f791e618474df3678176f73d707b9/e2e/$newMd2"
Set-RowHyperlink $wsZh '$C$3' $newXlf1Zh "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70a010d6535a8f2af430b4e65d2a91828e4ac020/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newXlf1Zh"
Remove-RowHyperlink $wsZh '$E$3'
Remove-RowHyperlink $wsZh '$F$3'
$wsZh.Range("E3").Clear()
$wsZh.Range("F3").Clear()

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$newXlf1De = "e04a53de-7ccd-4d91-8db5-1ad4f9b4db22.ce213699e41debf6b1b6a57a302349c00d6d9f03.de-de.xlf"
$handoffTimeDe = "2016-03-08 17:06:42"

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("C2").Value = $newXlf1De
$wsDe.Range("D2").Value = $handoffTimeDe
$wsDe.Range("G2").Value = $notHandedBack
Set-RowHyperlink $wsDe '$A$2' $newMd1 "https://github.com/OpenLocalizationTest/oltest/blob/7cae2979785f791e618474df3678176f73d707b9/e2e/$newMd1"
Set-RowHyperlink $wsDe '$C$2' $newXlf1De "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59d825233ea81a9556b3cff53dec0df91d5227b1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlf1De"
Remove-RowHyperlink $wsDe '$E$2'
Remove-RowHyperlink $wsDe '$F$2'
$wsDe.Range("E2").Clear()
$wsDe.Range("F2").Clear()

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("B3").Value = $newStatus
$wsDe.Range("C3").Value = $newXlf1De
$wsDe.Range("D3").Value = $handoffTimeDe
$wsDe.Range("G3").Value = $notHandedBack
Set-RowHyperlink $wsDe '$A$3' $newMd2 "https://github.com/OpenLocalizationTest/oltest/blob/7cae2979785f791e618474df3678176f73d707b9/e2e/$newMd2"
Set-RowHyperlink $wsDe '$C$3' $newXlf1De "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/59d825233ea81a9556b3cff53dec0df91d5227b1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newXlf1De"
Remove-RowHyperlink $wsDe '$E$3'
Remove-RowHyperlink $wsDe '$F$3'
$wsDe.Range("E3").Clear()
$wsDe.Range("F3").Clear()

Write-Host "Report generated for handoff."
